$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22: new timesheet entry ---
$ws.Cells.Item(21, 1).Copy() | Out-Null
$ws.Cells.Item(22, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(22, 1).Value = 42705
$ws.Cells.Item(22, 2).Value = "DEV"
$ws.Cells.Item(22, 3).Value = "Mijn aanbodpagina aangepast, footermenu, userlogin "
$ws.Cells.Item(22, 4).Value = 6

# --- Row 23: new timesheet entry ---
$ws.Cells.Item(21, 1).Copy() | Out-Null
$ws.Cells.Item(23, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(23, 1).Value = 42706
$ws.Cells.Item(23, 2).Value = "DEV"
$ws.Cells.Item(23, 3).Value = "Homepage, messages, accountpagina, algmene contacpagina, test online deployment (5u!!!!)"
$ws.Cells.Item(23, 4).Value = 10

# --- Update the total-hours formula to include the new rows ---
$ws.Cells.Item(3, 6).Formula = "=SUM(D2:D23)"

# --- Recalculate so the cached formula value is refreshed ---
$excel.CalculateFull()

# --- Update the active selection on the sheet ---
$ws.Range("C29").Select() | Out-Null
